$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E data range to text so numeric-looking strings (e.g. "594.08", "0.160", "1.00")
# are preserved exactly as text, matching the source data which stores these as strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.516.18'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.626.91'
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '594.08'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').Value = '167.75'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.88%  '
$ws.Range('D9').Value = '2.626.82'
$ws.Range('E9').Value = '  -1.57%  '
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('D11').Value = '0.160'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '27.65'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('D15').Value = '3.117.55'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').Value = '67.500.62'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').Value = '2.629.16'
$ws.Range('E18').Value = '  -1.60%  '
$ws.Range('D19').Value = '12.01'
$ws.Range('D20').Value = '8.03'
$ws.Range('E20').Value = '  +2.59%  '
$ws.Range('D21').Value = '358.23'
$ws.Range('E21').Value = '  -1.52%  '
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('E24').Value = '  -3.56%  '
$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').Value = '10.36'
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '69.98'
$ws.Range('E27').Value = '  -1.06%  '
$ws.Range('D28').Value = '2.762.87'
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').Value = '546.69'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('E33').Value = '  -2.17%  '
$ws.Range('D34').Value = '1.89'
$ws.Range('E34').Value = '  -1.46%  '
$ws.Range('E35').Value = '  +5.18%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  -2.32%  '
$ws.Range('D38').Value = '157.96'
$ws.Range('E38').Value = '  +1.58%  '
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('D41').Value = '18.28'
$ws.Range('E41').Value = '  +2.00%  '
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('D43').Value = '5.21'
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('D45').Value = '2.43'
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').Value = '152.94'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('E49').Value = '  -1.28%  '
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('E51').Value = '  -0.76%  '

# Restore default cell style (remove the temporary text-format style) so the
# saved cells have no explicit style index, matching the original formatting.
$ws.Range("D2:E51").Style = "Normal"
